# Apply weekly update: insert a new data row at row 15 (pushing existing
# rows 15-83 down to 16-84) and populate it with the latest week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 15; this shifts rows 15..83 down
# to 16..84 and copies formatting (e.g. the date style on column D) from
# the row above, same as Excel's normal "Insert Copied Cells"/row insert
# behavior.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with this week's record.
$ws.Cells.Item(15, 1).Value = 10
$ws.Cells.Item(15, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(15, 3).Value = "La Araucanía"
$ws.Cells.Item(15, 4).Value = 44565
$ws.Cells.Item(15, 5).Value = 9
$ws.Cells.Item(15, 6).Value = 100112022
$ws.Cells.Item(15, 7).Value = "Arveja Verde"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 20
$ws.Cells.Item(15, 11).Value = 26000
$ws.Cells.Item(15, 12).Value = 26000
$ws.Cells.Item(15, 13).Value = 26000
$ws.Cells.Item(15, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(15, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(15, 16).Value = 1040
$ws.Cells.Item(15, 17).Value = 25
$ws.Cells.Item(15, 18).Value = "Hortaliza"
